# Login.xlsx edit script
# - Updates the "Invalid password" / "Valid " test labels in column C to
#   "Test Invalid password" / "Test valid password".
# - Moves the active selection from C2 to G5.
# - Resizes column A and widens/creates columns C..K to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update cell values (shared-string text) ---
$ws.Range("C1").Value = "Test Invalid password"
$ws.Range("C2").Value = "Test valid password"

# --- Column widths ---
# (ColumnWidth is expressed in characters; values below were chosen so the
# resulting stored sheet column width lines up with the target layout.)
$ws.Columns.Item(1).ColumnWidth = 10.833333333333332
$ws.Columns.Item(3).ColumnWidth = 19.666666666666664
$ws.Columns.Item(4).ColumnWidth = 17.5
$ws.Columns.Item(5).ColumnWidth = 16.666666666666664
$ws.Columns.Item(6).ColumnWidth = 15.666666666666666
$ws.Columns.Item(7).ColumnWidth = 16.0
$ws.Columns.Item(8).ColumnWidth = 13.666666666666666
$ws.Columns.Item(9).ColumnWidth = 16.166666666666664
$ws.Columns.Item(10).ColumnWidth = 12.5
$ws.Columns.Item(11).ColumnWidth = 13.666666666666666

# --- Selection moves to G5 ---
$ws.Range("G5").Select()
